# Doing Updates for Financials
# Updates the KBAL yearly financials with refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KBAL")

# Balance Sheet - Net Receivables
$ws.Range("D43").Value = 62300

# Balance Sheet - Total Current Assets
$ws.Range("D46").Value = 207900

# Balance Sheet - Property Plant and Equipment
$ws.Range("D48").Value = 169000

# Balance Sheet - Total Assets
$ws.Range("D54").Value = 331500

# Balance Sheet - Other Current Liabilities
$ws.Range("D59").Value = 123800

# Balance Sheet - Total Current Liabilities
$ws.Range("D60").Value = 122700

# Balance Sheet - Total Liabilities
$ws.Range("D66").Value = 138400

# Cash Flow Statement - Capital Expenditures (full row refresh)
$ws.Range("D91").Value = -21600
$ws.Range("E91").Value = -11800
$ws.Range("F91").Value = -15000
$ws.Range("G91").Value = -31700
$ws.Range("H91").Value = -32900
$ws.Range("I91").Value = -27600
$ws.Range("J91").Value = -26900
